$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quote the "xD_i/xD_o" family labels in column A (rows 17-26) ---
# These cells keep referring to the same concept, but the label text is now
# wrapped in double quotes (e.g. 2D_i -> "2D_i").
$ws.Range("A17").Value2 = '"2D_i"'
$ws.Range("A18").Value2 = '"2D+_i"'
$ws.Range("A19").Value2 = '"3D_i"'
$ws.Range("A20").Value2 = '"3D+_i"'
$ws.Range("A21").Value2 = '"1D_o"'
$ws.Range("A22").Value2 = '"1D+_o"'
$ws.Range("A23").Value2 = '"2D_o"'
$ws.Range("A24").Value2 = '"2D+_o"'
$ws.Range("A25").Value2 = '"3D_o"'
$ws.Range("A26").Value2 = '"3D+_o"'

# --- Drop the stray, empty V10 cell (column V was only used there) ---
# Clearing it (rather than just clearing contents) removes the cell/column
# from the sheet's used range, which shrinks dimension + row spans from
# column 22 (V) down to column 20 (T), matching the updated layout.
$ws.Range("V10").Clear()

# --- Selection now spans the whole (shrunk) used range ---
$ws.Range("A1:T26").Select()
